$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.120.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.772.23'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.82%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.87%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.65%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3760'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.84%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.23'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.41%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3392'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.182'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.97%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07372'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.72%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.007'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.65%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.09%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.356'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.777.94'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.976'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.33%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001079'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.41%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06655'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.75%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.69%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.497'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.70%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.89%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.132.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.97%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.34%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.423'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.00%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.480'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.489'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.87%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.59%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '151.30'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.14%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.980.14'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '132.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.061'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.50%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.920'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.60%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.01%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.648'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.32%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.341'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.91%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6747'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.39%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06274'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02312'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.35%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2165'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.668'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.19%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.234'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.36'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.005'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6278'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.52%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.827'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.58%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.098'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.40%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '128.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.25%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07140'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.92%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.26%  '
